$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.679.32'
$ws.Range("E2").Value = '  +2.50%  '
$ws.Range("D3").Value = '2.042.25'
$ws.Range("E3").Value = '  +7.34%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '245.99'
$ws.Range("E5").Value = '  -0.64%  '
$ws.Range("D6").Value = '0.662'
$ws.Range("E6").Value = '  -4.75%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '45.04'
$ws.Range("E8").Value = '  +3.91%  '
$ws.Range("D9").Value = '60.35'
$ws.Range("E9").Value = '  +5.98%  '
$ws.Range("D10").Value = '0.361'
$ws.Range("E10").Value = '  +0.70%  '
$ws.Range("D11").Value = '0.0718'
$ws.Range("E11").Value = '  -5.07%  '
$ws.Range("D12").Value = '0.0985'
$ws.Range("E12").Value = '  -0.02%  '
$ws.Range("D13").Value = '14.52'
$ws.Range("E13").Value = '  -0.82%  '
$ws.Range("D14").Value = '2.332.14'
$ws.Range("E14").Value = '  +7.33%  '
$ws.Range("D15").Value = '0.807'
$ws.Range("E15").Value = '  +1.44%  '
$ws.Range("D16").Value = '2.031.60'
$ws.Range("E16").Value = '  +7.74%  '
$ws.Range("D17").Value = '4.88'
$ws.Range("E17").Value = '  -3.51%  '
$ws.Range("D18").Value = '36.542.90'
$ws.Range("E18").Value = '  +2.31%  '
$ws.Range("D19").Value = '71.10'
$ws.Range("E19").Value = '  -3.50%  '
$ws.Range("D20").Value = '0.0₃0812'
$ws.Range("E20").Value = '  -2.57%  '
$ws.Range("D21").Value = '236.51'
$ws.Range("E21").Value = '  -4.46%  '
$ws.Range("D22").Value = '12.57'
$ws.Range("E22").Value = '  -4.21%  '
$ws.Range("D23").Value = '4.89'
$ws.Range("E23").Value = '  -6.20%  '
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("D25").Value = '2.45'
$ws.Range("E25").Value = '  -9.11%  '
$ws.Range("D26").Value = '168.88'
$ws.Range("E26").Value = '  +1.18%  '
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").Value = '8.75'
$ws.Range("E27").Value = '  +0.15%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '19.83'
$ws.Range("E28").Value = '  +7.49%  '
$ws.Range("D29").Value = '1.94'
$ws.Range("E29").Value = '  -9.95%  '
$ws.Range("E30").Value = '  -5.60%  '
$ws.Range("D31").Value = '21.66'
$ws.Range("E31").Value = '  +50.21%  '
$ws.Range("D32").Value = '4.35'
$ws.Range("E32").Value = '  -1.89%  '
$ws.Range("D33").Value = '0.0579'
$ws.Range("E33").Value = '  -5.17%  '
$ws.Range("B34").Value = 'BinanceUSD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").Value = '  -0.08%  '
$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D35").Value = '0.0878'
$ws.Range("E35").Value = '  +18.63%  '
$ws.Range("E36").Value = '  +0.56%  '
$ws.Range("B37").Value = 'LidoDAOToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D37").Value = '2.22'
$ws.Range("E37").Value = '  +12.99%  '
$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").Value = '3.97'
$ws.Range("E38").Value = '  -7.17%  '
$ws.Range("D39").Value = '0.860'
$ws.Range("E39").Value = '  +0.06%  '
$ws.Range("D40").Value = '1.32'
$ws.Range("E40").Value = '  -11.50%  '
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = '0.0214'
$ws.Range("E41").Value = '  -6.78%  '
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Value = '95.93'
$ws.Range("E42").Value = '  -3.74%  '
$ws.Range("D43").Value = '1.11'
$ws.Range("E43").Value = '  +1.88%  '
$ws.Range("E44").Value = '  +15.63%  '
$ws.Range("D45").Value = '15.85'
$ws.Range("E45").Value = '  -7.25%  '
$ws.Range("D46").Value = '1.312.88'
$ws.Range("E46").Value = '  -0.34%  '
$ws.Range("D47").Value = '0.0813'
$ws.Range("E47").Value = '  -0.03%  '
$ws.Range("D48").Value = '2.80'
$ws.Range("E48").Value = '  +1.67%  '
$ws.Range("D49").Value = '2.226.47'
$ws.Range("E49").Value = '  +7.41%  '
$ws.Range("D51").Value = '3.82'
$ws.Range("E51").Value = '  +14.11%  '
